$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-valued cells need explicit Text number format first,
# otherwise Excel auto-converts strings like "87%" into the number 0.87.
$percentCells = @("H2", "H3", "H4", "H6", "H19", "H25", "H26", "H30", "H38", "H41")
foreach ($pc in $percentCells) {
    $ws.Range($pc).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-08 19:18:25"
$ws.Range("H2").Value = "87%"
$ws.Range("I2").Value = "4.5 mm"
$ws.Range("E3").Value = "2026-02-08 19:18:27"
$ws.Range("H3").Value = "96%"
$ws.Range("E4").Value = "2026-02-08 19:18:30"
$ws.Range("H4").Value = "67%"
$ws.Range("I4").Value = "1.7 mm"
$ws.Range("J4").Value = "1001.9 hPa"
$ws.Range("O4").Value = "10.4 °C"
$ws.Range("E5").Value = "2026-02-08 19:18:32"
$ws.Range("E6").Value = "2026-02-08 19:18:34"
$ws.Range("H6").Value = "65%"
$ws.Range("J6").Value = "1001.8 hPa"
$ws.Range("E7").Value = "2026-02-08 19:18:37"
$ws.Range("J7").Value = "1002.1 hPa"
$ws.Range("E8").Value = "2026-02-08 19:18:39"
$ws.Range("J8").Value = "1002.0 hPa"
$ws.Range("E9").Value = "2026-02-08 19:18:42"
$ws.Range("O9").Value = "10.4 °C"
$ws.Range("E10").Value = "2026-02-08 19:18:44"
$ws.Range("E11").Value = "2026-02-08 19:18:47"
$ws.Range("E12").Value = "2026-02-08 19:18:49"
$ws.Range("E13").Value = "2026-02-08 19:18:51"
$ws.Range("J13").Value = "1003.5 hPa"
$ws.Range("E14").Value = "2026-02-08 19:18:54"
$ws.Range("E15").Value = "2026-02-08 19:18:56"
$ws.Range("I15").Value = "0.1 mm"
$ws.Range("E16").Value = "2026-02-08 19:18:58"
$ws.Range("I16").Value = "3.2 mm"
$ws.Range("E17").Value = "2026-02-08 19:19:01"
$ws.Range("E18").Value = "2026-02-08 19:19:03"
$ws.Range("J18").Value = "1002.1 hPa"
$ws.Range("E19").Value = "2026-02-08 19:19:05"
$ws.Range("H19").Value = "89%"
$ws.Range("I19").Value = "10.3 mm"
$ws.Range("E20").Value = "2026-02-08 19:19:08"
$ws.Range("I20").Value = "8.3 mm"
$ws.Range("L20").Value = "54.4 km/h - 326º 18:43 TU"
$ws.Range("E21").Value = "2026-02-08 19:19:10"
$ws.Range("O21").Value = "5.5 °C"
$ws.Range("E22").Value = "2026-02-08 19:19:13"
$ws.Range("E23").Value = "2026-02-08 19:19:15"
$ws.Range("I23").Value = "4.7 mm"
$ws.Range("E24").Value = "2026-02-08 19:19:18"
$ws.Range("J24").Value = "1003.4 hPa"
$ws.Range("E25").Value = "2026-02-08 19:19:20"
$ws.Range("H25").Value = "78%"
$ws.Range("E26").Value = "2026-02-08 19:19:22"
$ws.Range("H26").Value = "68%"
$ws.Range("J26").Value = "1001.1 hPa"
$ws.Range("E27").Value = "2026-02-08 19:19:25"
$ws.Range("E28").Value = "2026-02-08 19:19:27"
$ws.Range("J28").Value = "1001.8 hPa"
$ws.Range("E29").Value = "2026-02-08 19:19:30"
$ws.Range("E30").Value = "2026-02-08 19:19:32"
$ws.Range("H30").Value = "69%"
$ws.Range("J30").Value = "1002.2 hPa"
$ws.Range("K30").Value = "10.3 MJ/m2"
$ws.Range("O30").Value = "10.0 °C"
$ws.Range("E31").Value = "2026-02-08 19:19:34"
$ws.Range("I31").Value = "0.1 mm"
$ws.Range("J31").Value = "1001.2 hPa"
$ws.Range("N31").Value = "8.1 °C 18:59 TU"
$ws.Range("O31").Value = "9.8 °C"
$ws.Range("E32").Value = "2026-02-08 19:19:37"
$ws.Range("E33").Value = "2026-02-08 19:19:39"
$ws.Range("E34").Value = "2026-02-08 19:19:41"
$ws.Range("E35").Value = "2026-02-08 19:19:44"
$ws.Range("J35").Value = "1004.2 hPa"
$ws.Range("E36").Value = "2026-02-08 19:19:46"
$ws.Range("J36").Value = "1002.2 hPa"
$ws.Range("E37").Value = "2026-02-08 19:19:49"
$ws.Range("J37").Value = "1003.0 hPa"
$ws.Range("E38").Value = "2026-02-08 19:19:51"
$ws.Range("H38").Value = "75%"
$ws.Range("I38").Value = "2.2 mm"
$ws.Range("E39").Value = "2026-02-08 19:19:54"
$ws.Range("E40").Value = "2026-02-08 19:19:56"
$ws.Range("J40").Value = "1003.5 hPa"
$ws.Range("E41").Value = "2026-02-08 19:19:58"
$ws.Range("H41").Value = "70%"
$ws.Range("J41").Value = "1002.3 hPa"
$ws.Range("E42").Value = "2026-02-08 19:20:01"
$ws.Range("E43").Value = "2026-02-08 19:20:03"
$ws.Range("E44").Value = "2026-02-08 19:20:05"
$ws.Range("I44").Value = "2.2 mm"
$ws.Range("E45").Value = "2026-02-08 19:20:08"
$ws.Range("J45").Value = "1004.3 hPa"
$ws.Range("E46").Value = "2026-02-08 19:20:11"
$ws.Range("J46").Value = "1003.9 hPa"
$ws.Range("O46").Value = "9.6 °C"
